$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A neutral, never-touched cell whose default (unstyled) look we borrow
# after forcing text-typed values below, so we don't leave a stray
# "quote prefix / text format" style on cells that should stay default.
$blankStyle = $ws.Range("Z100").Style

function Set-TextValue($range) {
    # Re-apply the default style so forcing text via a leading
    # apostrophe doesn't leave a lingering number-format style.
    $range.Style = $blankStyle
}

# Row 1: header labels replaced by sequential numbers 0-12
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
$ws.Range("M1").Value = 12

# Row 2: clear A2 and M2, set E2 to "Drive"
$ws.Range("A2").Value = ""
$ws.Range("E2").Value = "Drive"
$ws.Range("M2").Value = ""

# Row 3: becomes the old header-row labels (shifted down from row 1); clear J3 and M3
$ws.Range("A3").Value = "Lg."
$ws.Range("B3").Value = "Threading"
$ws.Range("C3").Value = "HeadDia."
$ws.Range("D3").Value = "Head Ht."
$ws.Range("E3").Value = "Style"
$ws.Range("F3").Value = "Size"
$ws.Range("G3").Value = "Tensile Strength, psi"
$ws.Range("H3").Value = "Specifications Met"
$ws.Range("I3").Value = "Pkg.Qty."
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "Pkg."
$ws.Range("M3").Value = ""

# Row 4
$ws.Range("A4").Value = "1/8"""
$ws.Range("J4").Value = "96710A050"
$ws.Range("K4").Value = "'`$14.89"
Set-TextValue $ws.Range("K4")
$ws.Range("L4").Value = "00-90"

# Row 5
$ws.Range("A5").Value = "3/16"""
$ws.Range("J5").Value = "96710A054"
$ws.Range("K5").Value = "'15.00"
Set-TextValue $ws.Range("K5")
$ws.Range("L5").Value = "00-90"

# Row 6: filled in completely (previously mostly blank "0-80" row)
$ws.Range("A6").Value = "1/4"""
$ws.Range("B6").Value = "Fully Threaded"
$ws.Range("C6").Value = "0.090"""
$ws.Range("D6").Value = "0.036"""
$ws.Range("E6").Value = "Torx Plus"
$ws.Range("F6").Value = "IP2, T2"
$ws.Range("G6").Value = "'70,000"
Set-TextValue $ws.Range("G6")
$ws.Range("H6").Value = "__"
$ws.Range("I6").Value = "'25"
Set-TextValue $ws.Range("I6")
$ws.Range("J6").Value = "96710A058"
$ws.Range("K6").Value = "'15.11"
Set-TextValue $ws.Range("K6")
$ws.Range("L6").Value = "00-90"

# Rows 7-10: set L column to "0-80"
$ws.Range("L7").Value = "0-80"
$ws.Range("L8").Value = "0-80"
$ws.Range("L9").Value = "0-80"
$ws.Range("L10").Value = "0-80"
